$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix trailing space in existing text "Atualizar diagramas UML "
for ($r = 1; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "Atualizar diagramas UML ") {
        $cell.Value = "Atualizar diagramas UML"
    }
}

# Add new row 32 data (row 31 was the last populated row before this edit)
$ws.Cells.Item(32, 1).Value = "implementar os casos de uso selecionados para a iteração"
$ws.Cells.Item(32, 2).Value = "Alta"
$ws.Cells.Item(32, 3).Value = 70
$ws.Cells.Item(32, 4).Value = "Iniciado"
$ws.Cells.Item(32, 5).Value = "E2"
$ws.Cells.Item(32, 6).Value = "Tarcísio/ Diógenes"
$ws.Cells.Item(32, 7).Value = 60
$ws.Cells.Item(32, 8).Value = 15

# Move selection to A33
$ws.Range("A33").Select()

# Re-apply the autofilter, which (as in the source repo's LibreOffice-driven
# history) stamps a new incrementally-numbered _FilterDatabase defined name
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0", "='Lista de Itens de Trabalho'!`$A`$1:`$I`$25")
